# "Plantilla de Ítems de Configuración" — actualización de reglas de nombrado.
# Reemplaza los nombres de ejemplo genéricos de la columna B (Reglas de
# Nombrado) por los nuevos nombres con marcadores entre "<>", agrega la
# regla de Papers que faltaba (B6) y añade cuatro filas con aclaraciones
# sobre los marcadores al final de la hoja.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Columna B: nuevas reglas de nombrado -------------------------------
$ws.Range("B2").Value  = "<NombreDelTema>.pdf"
$ws.Range("B3").Value  = "<NombreDelTema> - <Autor>.jpg"
$ws.Range("B4").Value  = "Trabajo_conceptual_N°_<Nro>.pdf"
$ws.Range("B5").Value  = "<NombreDelTrabajo>.pdf"
$ws.Range("B6").Value  = "<NombrePapers>.pdf"
$ws.Range("B7").Value  = "<NombreDelLibro> - <Autor> - <Edición>.<extension>"
$ws.Range("B8").Value  = "<NombreDelTema>.pdf"
$ws.Range("B9").Value  = "<NombreDelTema> - <Autor>.jpg"
$ws.Range("B10").Value = "Trabajo_practico_N°_<Nro>.pdf"
$ws.Range("B11").Value = "Estructura_de_Repositorio_V_<Version>.jpg"
$ws.Range("B12").Value = "Plantilla_de_items_de_Configuracion.xlsx"
$ws.Range("B13").Value = "NombreClase.<extension>"
$ws.Range("B14").Value = "US_Nro_<NombreUS>.docx"

# --- Nuevas filas al pie con la referencia de los marcadores ------------
$ws.Range("A19").Value = "<extension>: Indica la extension del archivo identificado como elemento de configuracion"
$ws.Range("A20").Value = "<Nro>: Indica el numero asignado a dicho elemento de configuracion"
$ws.Range("A21").Value = "<Version>: Indica la version del elemento de configuracion. Ejemplo: 1_00"
$ws.Range("A22").Value = "<NombreUS> Nombre de la User Story"

# --- Vista: zoom al 85% y selección en A16 -------------------------------
$excel.ActiveWindow.Zoom = 85
$ws.Range("A16").Select()
